# ActivityData_May22.xlsx - "Excel Test Driver and Model changes"
#
# Roll the test row's TransactionDate (C2) and PostingDate (D2) forward a
# month: 4/30/2022 -> 5/31/2022 (Excel serial dates 44681 -> 44712), and
# move the active selection from F2 to D2 to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TransactionActivity")

$ws.Range("C2").Value = 44712
$ws.Range("D2").Value = 44712

$ws.Range("D2").Select()
